$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J2").Value = 0.5
$ws.Range("J9").Value = 0.5
$ws.Range("J12").Value = -0.5
$ws.Range("J13").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("J16").Value = 0.5
$ws.Range("J22").Value = 0.5
